# Updated bootstrap Cap to 100V
# Insert a new BOM row (C5, C6 - 220N 100V capacitor) at row 10, pushing the
# rest of the Bill of Materials table down by one row, and fix up the
# print/filter ranges, autofilter, selection and manual page break that
# depend on the table's extents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10 (the Schottky diode row),
# shifting rows 10:53 down to 11:54.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the bootstrap capacitor entry.
# Columns are written B, C, A, E (then the numeric columns) to mirror the
# original authoring order.
$ws.Range("B10").Value = "C5, C6"
$ws.Range("C10").Value = "220N"
$ws.Range("A10").Value = "Unpolarized capacitor, 100V"
$ws.Range("E10").Value = "587-5013-1-ND"
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 0.26
$ws.Range("H10").Formula = "=G10*F10"

# Fix up the defined names that describe the table extents.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$6:`$H`$54"
    }
    if ($n.Name -eq "Sheet1!Print_Area") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$H`$23"
    }
}

# Re-apply the AutoFilter over the new, larger range.
$ws.AutoFilterMode = $false
$ws.Range("A6:H54").AutoFilter() | Out-Null

# Move the manual page break down by one row (was before row 28, now before row 29).
$ws.ResetAllPageBreaks()
$ws.Rows(29).PageBreak = 1

# Restore the active selection as recorded after the edit.
$ws.Range("H12").Select()
